# Apply weekly price-sheet update for "Hortaliza, Terminal La Palmera de La Serena - Melón"
# Rows 12-42 get refreshed figures (new week shifted in); rows 43-44 are newly added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = 44536
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 3200
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("N12").Value = '$/unidad'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 1450
$ws.Range("Q12").Value = 1

# Row 13
$ws.Range("D13").Value = 44536
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 2800
$ws.Range("K13").Value = 1200
$ws.Range("L13").Value = 1300
$ws.Range("M13").Value = 1250
$ws.Range("N13").Value = '$/unidad'
$ws.Range("O13").Value = 'Provincia de Limarí'
$ws.Range("P13").Value = 1250
$ws.Range("Q13").Value = 1

# Row 14
$ws.Range("D14").Value = 44515
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 26000
$ws.Range("M14").Value = 25500
$ws.Range("N14").Value = '$/caja 18 unidades'
$ws.Range("O14").Value = 'Provincia de Copiapó'
$ws.Range("P14").Value = 1417
$ws.Range("Q14").Value = 18

# Row 15
$ws.Range("D15").Value = 44246
$ws.Range("I15").Value = 'Extra'
$ws.Range("J15").Value = 5200
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1100
$ws.Range("M15").Value = 1050
$ws.Range("N15").Value = '$/unidad'
$ws.Range("O15").Value = 'Región de O''Higgins'
$ws.Range("P15").Value = 1050
$ws.Range("Q15").Value = 1

# Row 16
$ws.Range("D16").Value = 44246
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 4200
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = 875
$ws.Range("N16").Value = '$/unidad'
$ws.Range("O16").Value = 'Región de O''Higgins'
$ws.Range("P16").Value = 875
$ws.Range("Q16").Value = 1

# Row 17
$ws.Range("D17").Value = 44246
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 3200
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 775
$ws.Range("N17").Value = '$/unidad'
$ws.Range("O17").Value = 'Región de O''Higgins'
$ws.Range("P17").Value = 775
$ws.Range("Q17").Value = 1

# Row 18
$ws.Range("D18").Value = 44225
$ws.Range("I18").Value = 'Extra'
$ws.Range("J18").Value = 4000
$ws.Range("K18").Value = 850
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 875
$ws.Range("N18").Value = '$/unidad'
$ws.Range("O18").Value = 'Región de O''Higgins'
$ws.Range("P18").Value = 875
$ws.Range("Q18").Value = 1

# Row 19
$ws.Range("D19").Value = 44225
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 3200
$ws.Range("K19").Value = 650
$ws.Range("L19").Value = 700
$ws.Range("M19").Value = 675
$ws.Range("N19").Value = '$/unidad'
$ws.Range("O19").Value = 'Región de O''Higgins'
$ws.Range("P19").Value = 675
$ws.Range("Q19").Value = 1

# Row 20
$ws.Range("D20").Value = 44225
$ws.Range("I20").Value = 'Segunda'
$ws.Range("J20").Value = 2600
$ws.Range("K20").Value = 450
$ws.Range("L20").Value = 500
$ws.Range("M20").Value = 475
$ws.Range("N20").Value = '$/unidad'
$ws.Range("O20").Value = 'Región de O''Higgins'
$ws.Range("P20").Value = 475
$ws.Range("Q20").Value = 1

# Row 21
$ws.Range("D21").Value = 44232
$ws.Range("I21").Value = 'Extra'
$ws.Range("J21").Value = 6000
$ws.Range("K21").Value = 950
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 975
$ws.Range("N21").Value = '$/unidad'
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 975
$ws.Range("Q21").Value = 1

# Row 22
$ws.Range("D22").Value = 44232
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 875
$ws.Range("N22").Value = '$/unidad'
$ws.Range("O22").Value = 'Región Metropolitana'
$ws.Range("P22").Value = 875
$ws.Range("Q22").Value = 1

# Row 23
$ws.Range("D23").Value = 44232
$ws.Range("I23").Value = 'Segunda'
$ws.Range("J23").Value = 4000
$ws.Range("K23").Value = 750
$ws.Range("L23").Value = 800
$ws.Range("M23").Value = 775
$ws.Range("N23").Value = '$/unidad'
$ws.Range("O23").Value = 'Región Metropolitana'
$ws.Range("P23").Value = 775
$ws.Range("Q23").Value = 1

# Row 24
$ws.Range("D24").Value = 44242
$ws.Range("I24").Value = 'Extra'
$ws.Range("J24").Value = 5200
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1100
$ws.Range("M24").Value = 1050
$ws.Range("N24").Value = '$/unidad'
$ws.Range("O24").Value = 'Región Metropolitana'
$ws.Range("P24").Value = 1050
$ws.Range("Q24").Value = 1

# Row 25
$ws.Range("D25").Value = 44242
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 4000
$ws.Range("K25").Value = 850
$ws.Range("L25").Value = 900
$ws.Range("M25").Value = 875
$ws.Range("N25").Value = '$/unidad'
$ws.Range("O25").Value = 'Región Metropolitana'
$ws.Range("P25").Value = 875
$ws.Range("Q25").Value = 1

# Row 26
$ws.Range("D26").Value = 44242
$ws.Range("I26").Value = 'Segunda'
$ws.Range("J26").Value = 3600
$ws.Range("K26").Value = 750
$ws.Range("L26").Value = 800
$ws.Range("M26").Value = 775
$ws.Range("N26").Value = '$/unidad'
$ws.Range("O26").Value = 'Región Metropolitana'
$ws.Range("P26").Value = 775
$ws.Range("Q26").Value = 1

# Row 27
$ws.Range("D27").Value = 44204
$ws.Range("I27").Value = 'Extra'
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1100
$ws.Range("M27").Value = 1050
$ws.Range("N27").Value = '$/unidad'
$ws.Range("O27").Value = 'Región de O''Higgins'
$ws.Range("P27").Value = 1050
$ws.Range("Q27").Value = 1

# Row 28
$ws.Range("D28").Value = 44204
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 800
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = 850
$ws.Range("N28").Value = '$/unidad'
$ws.Range("O28").Value = 'Región de O''Higgins'
$ws.Range("P28").Value = 850
$ws.Range("Q28").Value = 1

# Row 29
$ws.Range("D29").Value = 44204
$ws.Range("I29").Value = 'Segunda'
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 700
$ws.Range("L29").Value = 750
$ws.Range("M29").Value = 725
$ws.Range("N29").Value = '$/unidad'
$ws.Range("O29").Value = 'Región de O''Higgins'
$ws.Range("P29").Value = 725
$ws.Range("Q29").Value = 1

# Row 30
$ws.Range("D30").Value = 44204
$ws.Range("I30").Value = 'Super'
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 1200
$ws.Range("L30").Value = 1300
$ws.Range("M30").Value = 1250
$ws.Range("N30").Value = '$/unidad'
$ws.Range("O30").Value = 'Región de O''Higgins'
$ws.Range("P30").Value = 1250
$ws.Range("Q30").Value = 1

# Row 31
$ws.Range("D31").Value = 44200
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 750
$ws.Range("L31").Value = 800
$ws.Range("M31").Value = 775
$ws.Range("N31").Value = '$/unidad'
$ws.Range("O31").Value = 'Región de O''Higgins'
$ws.Range("P31").Value = 775
$ws.Range("Q31").Value = 1

# Row 32
$ws.Range("D32").Value = 44200
$ws.Range("I32").Value = 'Segunda'
$ws.Range("J32").Value = 1600
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 650
$ws.Range("M32").Value = 625
$ws.Range("N32").Value = '$/unidad'
$ws.Range("O32").Value = 'Región de O''Higgins'
$ws.Range("P32").Value = 625
$ws.Range("Q32").Value = 1

# Row 33
$ws.Range("D33").Value = 44522
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 23000
$ws.Range("L33").Value = 24000
$ws.Range("M33").Value = 23500
$ws.Range("N33").Value = '$/caja 16 unidades'
$ws.Range("O33").Value = 'Provincia de Copiapó'
$ws.Range("P33").Value = 1469
$ws.Range("Q33").Value = 16

# Row 34
$ws.Range("D34").Value = 44218
$ws.Range("I34").Value = 'Extra'
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = 950
$ws.Range("N34").Value = '$/unidad'
$ws.Range("O34").Value = 'Región de O''Higgins'
$ws.Range("P34").Value = 950
$ws.Range("Q34").Value = 1

# Row 35
$ws.Range("D35").Value = 44218
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 700
$ws.Range("L35").Value = 800
$ws.Range("M35").Value = 750
$ws.Range("N35").Value = '$/unidad'
$ws.Range("O35").Value = 'Región de O''Higgins'
$ws.Range("P35").Value = 750
$ws.Range("Q35").Value = 1

# Row 36
$ws.Range("D36").Value = 44218
$ws.Range("I36").Value = 'Segunda'
$ws.Range("J36").Value = 2400
$ws.Range("K36").Value = 550
$ws.Range("L36").Value = 600
$ws.Range("M36").Value = 575
$ws.Range("N36").Value = '$/unidad'
$ws.Range("O36").Value = 'Región de O''Higgins'
$ws.Range("P36").Value = 575
$ws.Range("Q36").Value = 1

# Row 37
$ws.Range("D37").Value = 44525
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 23000
$ws.Range("L37").Value = 24000
$ws.Range("M37").Value = 23500
$ws.Range("N37").Value = '$/caja 16 unidades'
$ws.Range("O37").Value = 'Provincia de Copiapó'
$ws.Range("P37").Value = 1469
$ws.Range("Q37").Value = 16

# Row 38
$ws.Range("D38").Value = 44525
$ws.Range("I38").Value = 'Segunda'
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 19000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = 19500
$ws.Range("N38").Value = '$/caja 24 unidades'
$ws.Range("O38").Value = 'Provincia de Copiapó'
$ws.Range("P38").Value = 812
$ws.Range("Q38").Value = 24

# Row 39
$ws.Range("D39").Value = 44239
$ws.Range("I39").Value = 'Extra'
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 1000
$ws.Range("L39").Value = 1100
$ws.Range("M39").Value = 1050
$ws.Range("N39").Value = '$/unidad'
$ws.Range("O39").Value = 'Región Metropolitana'
$ws.Range("P39").Value = 1050
$ws.Range("Q39").Value = 1

# Row 40
$ws.Range("D40").Value = 44239
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 850
$ws.Range("L40").Value = 900
$ws.Range("M40").Value = 875
$ws.Range("N40").Value = '$/unidad'
$ws.Range("O40").Value = 'Región Metropolitana'
$ws.Range("P40").Value = 875
$ws.Range("Q40").Value = 1

# Row 41
$ws.Range("D41").Value = 44239
$ws.Range("I41").Value = 'Segunda'
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 750
$ws.Range("L41").Value = 800
$ws.Range("M41").Value = 775
$ws.Range("N41").Value = '$/unidad'
$ws.Range("O41").Value = 'Región Metropolitana'
$ws.Range("P41").Value = 775
$ws.Range("Q41").Value = 1

# Row 42
$ws.Range("D42").Value = 44211
$ws.Range("I42").Value = 'Extra'
$ws.Range("J42").Value = 4000
$ws.Range("K42").Value = 900
$ws.Range("L42").Value = 1000
$ws.Range("M42").Value = 950
$ws.Range("N42").Value = '$/unidad'
$ws.Range("O42").Value = 'Región de O''Higgins'
$ws.Range("P42").Value = 950
$ws.Range("Q42").Value = 1

# Row 43 (new row)
$ws.Range("A43").Value = 8
$ws.Range("B43").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C43").Value = 'Coquimbo'
$ws.Range("D43").Value = 44211
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 100112027
$ws.Range("G43").Value = 'Melón'
$ws.Range("H43").Value = 'Tuna'
$ws.Range("I43").Value = 'Primera'
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 700
$ws.Range("L43").Value = 800
$ws.Range("M43").Value = 750
$ws.Range("N43").Value = '$/unidad'
$ws.Range("O43").Value = 'Región de O''Higgins'
$ws.Range("P43").Value = 750
$ws.Range("Q43").Value = 1
$ws.Range("R43").Value = 'Hortaliza'

# Row 44 (new row)
$ws.Range("A44").Value = 8
$ws.Range("B44").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C44").Value = 'Coquimbo'
$ws.Range("D44").Value = 44211
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E44").Value = 4
$ws.Range("F44").Value = 100112027
$ws.Range("G44").Value = 'Melón'
$ws.Range("H44").Value = 'Tuna'
$ws.Range("I44").Value = 'Segunda'
$ws.Range("J44").Value = 3000
$ws.Range("K44").Value = 550
$ws.Range("L44").Value = 600
$ws.Range("M44").Value = 575
$ws.Range("N44").Value = '$/unidad'
$ws.Range("O44").Value = 'Región de O''Higgins'
$ws.Range("P44").Value = 575
$ws.Range("Q44").Value = 1
$ws.Range("R44").Value = 'Hortaliza'

